# Automatische test-sync: 2025-06-24 21:35:50
# Adds a new log entry (row 33) to the "Logs" sheet, extends the
# conditional-formatting ranges to cover the new row, and updates the
# corresponding aggregate count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# New row of data in the Logs sheet (row 33).
$logs.Range("A33").Value = "Ruilen van product"
$logs.Range("B33").Value = "mailmind.test@zohomail.eu"
$logs.Range("C33").Value = "Kan ik dit product ruilen voor een andere maat?"
$logs.Range("D33").Value = "Retour / Terugbetaling"
$logs.Range("F33").Value = "2025-06-24 21:34:51"
$logs.Range("G33").Value = "Nee"

# Extend the existing conditional formatting rules so they keep covering
# column D and G down to the newly added row 33 (previously D2:D32 / G2:G32).
$dFormatConditions = $logs.Range("D2:D32").FormatConditions
for ($i = 1; $i -le $dFormatConditions.Count; $i++) {
    $dFormatConditions.Item($i).ModifyAppliesToRange($logs.Range("D2:D33"))
}

$gFormatConditions = $logs.Range("G2:G32").FormatConditions
for ($i = 1; $i -le $gFormatConditions.Count; $i++) {
    $gFormatConditions.Item($i).ModifyAppliesToRange($logs.Range("G2:G33"))
}

# Update the Dashboard aggregate count for "Retour / Terugbetaling" (14 -> 15).
$dashboard.Range("B2").Value = 15
